# enhancements in call_queue monitoring
# Append three new "minimal_interaction" call records for Vanshika panjwani
# to the Incomplete Calls sheet (rows 4-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Name, Phone Number, Address, Age, Gender, Call Timestamp,
#             Call Duration (seconds), Reason, Notes
$newRows = @(
    @("Vanshika panjwani", "917823844614", "24 MG Road, Bengaluru", "28", "Male", "2025-06-25 19:35:46", 0,  "minimal_interaction", "Very few exchanges in conversation"),
    @("Vanshika panjwani", "917823844614", "24 MG Road, Bengaluru", "28", "Male", "2025-06-25 20:32:23", 88, "minimal_interaction", "Very few exchanges in conversation"),
    @("Vanshika panjwani", "917823844614", "24 MG Road, Bengaluru", "28", "Male", "2025-06-25 20:34:46", 1,  "minimal_interaction", "Very few exchanges in conversation")
)

$row = 4
foreach ($entry in $newRows) {
    # Columns A, C, E, H, I are plain text.
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 8).Value = $entry[7]
    $ws.Cells.Item($row, 9).Value = $entry[8]

    # Columns B (phone number) and D (age) look numeric, but must be stored
    # as text (no leading "+", matching the source data). A leading
    # apostrophe forces Excel to keep them as text instead of auto-coercing
    # them into numbers.
    $ws.Cells.Item($row, 2).Value = "'" + $entry[1]
    $ws.Cells.Item($row, 4).Value = "'" + $entry[3]

    # Column F (call timestamp) stays as text automatically.
    $ws.Cells.Item($row, 6).Value = $entry[5]

    # Column G (call duration) is a real number.
    $ws.Cells.Item($row, 7).Value = $entry[6]

    $row++
}
